$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# For each data row, the KYC update now supports cases where the investor's
# name or PAN is nil: correspondence-address codes get an "x" suffix, and
# both "Verified" and "Update Only" are marked "Yes".
for ($r = 2; $r -le 5; $r++) {
    $addr = $ws.Cells.Item($r, 5).Value2
    $ws.Cells.Item($r, 5).Value = "$($addr)x"
    $ws.Cells.Item($r, 14).Value = "Yes"
    $ws.Cells.Item($r, 15).Value = "Yes"
}

$ws.Range("O5").Select()
